# Update cryptocurrency price (D) and 1h volume change (E) columns
# for rows 2-51 on the active worksheet, matching the source refresh.

function Set-TextValue {
    # Writes $text into $range as plain text, exactly as-is.
    # A leading apostrophe forces Excel to keep number-looking
    # strings (e.g. "1.00", "91.50", "0.0000102") as text instead
    # of silently converting them to numeric values.
    param($range, [string]$text)
    if ($text -match '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*%?\s*$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "46.496.65"
Set-TextValue $ws.Range("E2") "  +0.80%  "

Set-TextValue $ws.Range("D3") "2.611.97"
Set-TextValue $ws.Range("E3") "  +10.55%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "307.84"
Set-TextValue $ws.Range("E5") "  +2.13%  "

Set-TextValue $ws.Range("D6") "101.51"
Set-TextValue $ws.Range("E6") "  +1.59%  "

Set-TextValue $ws.Range("D7") "0.606"
Set-TextValue $ws.Range("E7") "  +6.35%  "

Set-TextValue $ws.Range("D8") "  "
Set-TextValue $ws.Range("E8") "  +0.02%  "

Set-TextValue $ws.Range("D9") "0.582"
Set-TextValue $ws.Range("E9") "  +13.05%  "

Set-TextValue $ws.Range("D10") "38.94"
Set-TextValue $ws.Range("E10") "  +12.67%  "

Set-TextValue $ws.Range("D11") "0.0843"
Set-TextValue $ws.Range("E11") "  +5.75%  "

Set-TextValue $ws.Range("D12") "8.29"
Set-TextValue $ws.Range("E12") "  +16.19%  "

Set-TextValue $ws.Range("D13") "3.012.63"
Set-TextValue $ws.Range("E13") "  +10.76%  "

Set-TextValue $ws.Range("D14") "  "
Set-TextValue $ws.Range("E14") "  +2.07%  "

Set-TextValue $ws.Range("D15") "2.620.03"
Set-TextValue $ws.Range("E15") "  +10.95%  "

Set-TextValue $ws.Range("D16") "0.908"
Set-TextValue $ws.Range("E16") "  +11.34%  "

Set-TextValue $ws.Range("D17") "14.94"
Set-TextValue $ws.Range("E17") "  +9.63%  "

Set-TextValue $ws.Range("D18") "46.634.95"
Set-TextValue $ws.Range("E18") "  +1.26%  "

Set-TextValue $ws.Range("D19") "13.34"
Set-TextValue $ws.Range("E19") "  +3.69%  "

Set-TextValue $ws.Range("D20") "0.0000102"
Set-TextValue $ws.Range("E20") "  +5.14%  "

Set-TextValue $ws.Range("D21") "6.75"
Set-TextValue $ws.Range("E21") "  +11.81%  "

Set-TextValue $ws.Range("D22") "71.49"
Set-TextValue $ws.Range("E22") "  +5.70%  "

Set-TextValue $ws.Range("D23") "257.85"
Set-TextValue $ws.Range("E23") "  +5.09%  "

Set-TextValue $ws.Range("D24") "3.03"
Set-TextValue $ws.Range("E24") "  +7.74%  "

Set-TextValue $ws.Range("D25") "  "
Set-TextValue $ws.Range("E25") "  +15.93%  "

Set-TextValue $ws.Range("D26") "28.37"
Set-TextValue $ws.Range("E26") "  +35.19%  "

Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  -0.05%  "

Set-TextValue $ws.Range("D28") "10.54"
Set-TextValue $ws.Range("E28") "  +7.47%  "

Set-TextValue $ws.Range("D29") "39.88"
Set-TextValue $ws.Range("E29") "  -0.20%  "

Set-TextValue $ws.Range("D30") "  "
Set-TextValue $ws.Range("E30") "  +3.74%  "

Set-TextValue $ws.Range("D31") "6.20"
Set-TextValue $ws.Range("E31") "  +12.29%  "

Set-TextValue $ws.Range("D32") "3.74"
Set-TextValue $ws.Range("E32") "  -1.15%  "

Set-TextValue $ws.Range("D33") "2.34"
Set-TextValue $ws.Range("E33") "  +23.28%  "

Set-TextValue $ws.Range("D34") "2.95"
Set-TextValue $ws.Range("E34") "  +4.85%  "

Set-TextValue $ws.Range("D36") "151.11"
Set-TextValue $ws.Range("E36") "  +3.37%  "

Set-TextValue $ws.Range("D37") "  "
Set-TextValue $ws.Range("E37") "  +4.64%  "

Set-TextValue $ws.Range("D38") "  "
Set-TextValue $ws.Range("E38") "  +4.88%  "

Set-TextValue $ws.Range("D39") "4.21"
Set-TextValue $ws.Range("E39") "  +7.00%  "

Set-TextValue $ws.Range("D40") "15.80"
Set-TextValue $ws.Range("E40") "  +4.92%  "

Set-TextValue $ws.Range("D41") "3.64"
Set-TextValue $ws.Range("E41") "  +13.16%  "

Set-TextValue $ws.Range("D43") "2.062.82"
Set-TextValue $ws.Range("E43") "  +8.09%  "

Set-TextValue $ws.Range("D44") "19.14"
Set-TextValue $ws.Range("E44") "  +31.71%  "

Set-TextValue $ws.Range("D45") "  "
Set-TextValue $ws.Range("E45") "  -0.03%  "

Set-TextValue $ws.Range("D46") "91.50"
Set-TextValue $ws.Range("E46") "  -1.17%  "

Set-TextValue $ws.Range("D47") "1.80"
Set-TextValue $ws.Range("E47") "  +0.22%  "

Set-TextValue $ws.Range("D48") "9.23"
Set-TextValue $ws.Range("E48") "  +10.48%  "

Set-TextValue $ws.Range("D49") "110.43"
Set-TextValue $ws.Range("E49") "  +12.82%  "

Set-TextValue $ws.Range("D50") "  "
Set-TextValue $ws.Range("E50") "  +8.45%  "

Set-TextValue $ws.Range("D51") "2.869.96"
Set-TextValue $ws.Range("E51") "  +10.75%  "
